# The commit inserts one new weekly price-record row into the "Papa"
# sheet at row 127 (pushing the existing rows 127-191 down to 128-192),
# and fills the new row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 127, shifting rows 127-191 down to
# 128-192 (Excel copies formatting from the row above automatically).
$ws.Rows.Item(127).EntireRow.Insert()

# Populate the newly inserted row 127 with the new record.
$ws.Cells.Item(127, 1).Value  = 7
$ws.Cells.Item(127, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(127, 3).Value  = "Ñuble"
$ws.Cells.Item(127, 4).Value  = 44455
$ws.Cells.Item(127, 5).Value  = 16
$ws.Cells.Item(127, 6).Value  = 100114001
$ws.Cells.Item(127, 7).Value  = "Papa"
$ws.Cells.Item(127, 8).Value  = "Patagonia"
$ws.Cells.Item(127, 9).Value  = "1a (guarda)"
$ws.Cells.Item(127, 10).Value = 600
$ws.Cells.Item(127, 11).Value = 7000
$ws.Cells.Item(127, 12).Value = 7500
$ws.Cells.Item(127, 13).Value = 7250
$ws.Cells.Item(127, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(127, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(127, 16).Value = 290
$ws.Cells.Item(127, 17).Value = 25
$ws.Cells.Item(127, 18).Value = "Hortaliza"
